# Apply scraper refresh (08:34:05 run, Línea 141 - 828) to all three sheets.
$wb = $excel.ActiveWorkbook

function Set-Row($ws, $r, $a, $b, $c, $d, $e) {
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
}

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2, 1).Value = "Última actualización: 08:34:05"
$ws.Cells.Item(3, 1).Value = "Total filas: 125"

Set-Row $ws 39 "06:21:22" "06:29" "23_HERNANDEZ" 8 "LP1912"
Set-Row $ws 40 "06:21:22" "06:29" "86_EST CHICA-ESC AGRARIA" 8 "LP1912"
Set-Row $ws 41 "05:52:07" "06:30" "86_EST CHICA-ESC AGRARIA" 38 "LP1912"
Set-Row $ws 42 "06:21:22" "06:31" "16_SANTA ANA" 10 "LP1912"
Set-Row $ws 43 "04:48:57" "06:43" "225_C ROCA-H SUR" 115 "LP1912"
Set-Row $ws 44 "06:21:22" "06:44" "225_C ROCA-H SUR" 23 "LP1912"
Set-Row $ws 45 "06:21:22" "06:46" "215C_EL PATO" 25 "LP1912"
Set-Row $ws 46 "05:52:07" "06:47" "215C_EL PATO" 55 "LP1912"
Set-Row $ws 47 "06:59:37" "06:59" "14_ABASTO" 0 "LP1912"
Set-Row $ws 48 "06:59:37" "06:59" "23_HERNANDEZ" 0 "LP1912"
Set-Row $ws 49 "05:52:07" "07:00" "14_ABASTO" 68 "LP1912"
Set-Row $ws 50 "06:49:33" "07:01" "16_SANTA ANA" 12 "LP1912"
Set-Row $ws 51 "06:49:33" "07:04" "23_HERNANDEZ" 15 "LP1912"
Set-Row $ws 52 "05:52:07" "07:05" "23_HERNANDEZ" 73 "LP1912"
Set-Row $ws 53 "06:59:37" "07:05" "15_ABASTO" 6 "LP1912"
Set-Row $ws 54 "06:59:37" "07:07" "225_GOMEZ" 8 "LP1912"
Set-Row $ws 55 "06:59:37" "07:11" "215A_EL PATO" 12 "LP1912"
Set-Row $ws 56 "05:52:07" "07:12" "215A_EL PATO" 80 "LP1912"
Set-Row $ws 57 "06:59:37" "07:15" "11_ETCHEVERRY" 16 "LP1912"
Set-Row $ws 58 "06:59:37" "07:16" "16_SANTA ANA" 17 "LP1912"
Set-Row $ws 59 "05:52:07" "07:16" "11_ETCHEVERRY" 84 "LP1912"
Set-Row $ws 60 "06:59:37" "07:21" "26_HERNANDEZ" 22 "LP1912"
Set-Row $ws 61 "06:59:37" "07:23" "10_OLMOS" 24 "LP1912"
Set-Row $ws 62 "07:28:14" "07:30" "11_ETCHEVERRY" 2 "LP1912"
Set-Row $ws 63 "06:59:37" "07:31" "16_SANTA ANA" 32 "LP1912"
Set-Row $ws 64 "06:59:37" "07:31" "11_ETCHEVERRY" 32 "LP1912"
Set-Row $ws 65 "07:28:14" "07:32" "84_COLONIA URQUIZA-ESC 49" 4 "LP1912"
Set-Row $ws 66 "05:52:07" "07:32" "11_ETCHEVERRY" 100 "LP1912"
Set-Row $ws 67 "05:52:07" "07:32" "16_SANTA ANA" 100 "LP1912"
Set-Row $ws 68 "07:28:14" "07:35" "23_HERNANDEZ" 7 "LP1912"
Set-Row $ws 69 "06:59:37" "07:36" "27_EL RETIRO" 37 "LP1912"
Set-Row $ws 70 "07:28:14" "07:37" "27_EL RETIRO" 9 "LP1912"
Set-Row $ws 71 "07:28:14" "07:39" "10_OLMOS" 11 "LP1912"
Set-Row $ws 72 "06:59:37" "07:47" "14_ABASTO" 48 "LP1912"
Set-Row $ws 73 "07:28:14" "07:47" "16_SANTA ANA" 19 "LP1912"
Set-Row $ws 74 "07:28:14" "07:48" "14_ABASTO" 20 "LP1912"
Set-Row $ws 75 "07:51:34" "07:51" "215D_EL PATO" 0 "LP1912"
Set-Row $ws 76 "07:51:34" "07:51" "10_OLMOS" 0 "LP1912"
Set-Row $ws 77 "07:28:14" "07:55" "10_OLMOS" 27 "LP1912"
Set-Row $ws 78 "07:28:14" "08:00" "23_HERNANDEZ" 32 "LP1912"
Set-Row $ws 79 "07:51:34" "08:01" "23_HERNANDEZ" 10 "LP1912"
Set-Row $ws 80 "07:51:34" "08:03" "11_ETCHEVERRY" 12 "LP1912"
Set-Row $ws 81 "06:59:37" "08:06" "23_HERNANDEZ" 67 "LP1912"
Set-Row $ws 82 "07:51:34" "08:10" "16_SANTA ANA" 19 "LP1912"
Set-Row $ws 83 "07:28:14" "08:11" "16_SANTA ANA" 43 "LP1912"
Set-Row $ws 84 "07:51:34" "08:12" "15_ABASTO" 21 "LP1912"
Set-Row $ws 85 "07:51:34" "08:13" "10_OLMOS" 22 "LP1912"
Set-Row $ws 86 "08:13:38" "08:21" "26_HERNANDEZ" 8 "LP1912"
Set-Row $ws 87 "08:13:38" "08:22" "16_P MOR-SANTA ANA" 9 "LP1912"
Set-Row $ws 88 "08:13:38" "08:23" "215B_EL PATO" 10 "LP1912"
Set-Row $ws 89 "07:28:14" "08:23" "16_P MOR-SANTA ANA" 55 "LP1912"
Set-Row $ws 90 "08:13:38" "08:27" "84_COLONIA URQUIZA-ESC 49" 14 "LP1912"
Set-Row $ws 91 "07:51:34" "08:30" "23_HERNANDEZ" 39 "LP1912"
Set-Row $ws 92 "08:13:38" "08:33" "10_OLMOS" 20 "LP1912"
Set-Row $ws 93 "08:13:38" "08:36" "23_HERNANDEZ" 23 "LP1912"
Set-Row $ws 94 "08:34:05" "08:42" "81_EL PELIGRO" 8 "LP1912"
Set-Row $ws 95 "08:13:38" "08:43" "14_ABASTO" 30 "LP1912"
Set-Row $ws 96 "08:34:05" "08:44" "14_ABASTO" 10 "LP1912"
Set-Row $ws 97 "08:34:05" "08:53" "10_OLMOS" 19 "LP1912"
Set-Row $ws 98 "08:34:05" "08:54" "17_ROMERO" 20 "LP1912"
Set-Row $ws 99 "08:13:38" "09:01" "23_HERNANDEZ" 48 "LP1912"
Set-Row $ws 100 "08:13:38" "09:01" "215A_EL PATO" 48 "LP1912"
Set-Row $ws 101 "08:34:05" "09:02" "215A_EL PATO" 28 "LP1912"
Set-Row $ws 102 "08:13:38" "09:03" "11_ETCHEVERRY" 50 "LP1912"
Set-Row $ws 103 "08:34:05" "09:04" "11_ETCHEVERRY" 30 "LP1912"
Set-Row $ws 104 "08:34:05" "09:05" "23_HERNANDEZ" 31 "LP1912"
Set-Row $ws 105 "08:13:38" "09:10" "16_P MOR-SANTA ANA" 57 "LP1912"
Set-Row $ws 106 "08:34:05" "09:11" "16_P MOR-SANTA ANA" 37 "LP1912"
Set-Row $ws 107 "08:34:05" "09:13" "10_OLMOS" 39 "LP1912"
Set-Row $ws 108 "08:13:38" "09:16" "27_EL RETIRO" 63 "LP1912"
Set-Row $ws 109 "08:34:05" "09:17" "27_EL RETIRO" 43 "LP1912"
Set-Row $ws 110 "08:34:05" "09:21" "26_HERNANDEZ" 47 "LP1912"
Set-Row $ws 111 "08:13:38" "09:22" "17_ROMERO" 69 "LP1912"
Set-Row $ws 112 "08:34:05" "09:23" "16_SANTA ANA" 49 "LP1912"
Set-Row $ws 113 "07:28:14" "09:23" "17_ROMERO" 115 "LP1912"
Set-Row $ws 114 "08:13:38" "09:23" "11_ETCHEVERRY" 70 "LP1912"
Set-Row $ws 115 "08:34:05" "09:24" "11_ETCHEVERRY" 50 "LP1912"
Set-Row $ws 116 "08:34:05" "09:32" "15_ABASTO" 58 "LP1912"
Set-Row $ws 117 "08:34:05" "09:33" "10_OLMOS" 59 "LP1912"
Set-Row $ws 118 "08:13:38" "09:34" "16_SANTA ANA" 81 "LP1912"
Set-Row $ws 119 "08:34:05" "09:35" "23_HERNANDEZ" 61 "LP1912"
Set-Row $ws 120 "08:34:05" "09:35" "16_SANTA ANA" 61 "LP1912"
Set-Row $ws 121 "08:13:38" "09:41" "215C_EL PATO" 88 "LP1912"
Set-Row $ws 122 "08:34:05" "09:42" "215C_EL PATO" 68 "LP1912"
Set-Row $ws 123 "08:13:38" "09:43" "14_ABASTO" 90 "LP1912"
Set-Row $ws 124 "08:34:05" "09:44" "14_ABASTO" 70 "LP1912"
Set-Row $ws 125 "08:13:38" "09:58" "16_SANTA ANA" 105 "LP1912"
Set-Row $ws 126 "08:34:05" "10:11" "16_P MOR-SANTA ANA" 97 "LP1912"
Set-Row $ws 127 "08:34:05" "10:12" "15_ABASTO" 98 "LP1912"
Set-Row $ws 128 "08:34:05" "10:21" "26_HERNANDEZ" 107 "LP1912"
Set-Row $ws 129 "08:34:05" "10:22" "17_ROMERO" 108 "LP1912"
Set-Row $ws 130 "08:34:05" "10:27" "215A_EL PATO" 113 "LP1912"

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2, 1).Value = "Última actualización: 08:34:05"
$ws.Cells.Item(3, 1).Value = "Total filas: 21"

Set-Row $ws 23 "08:34:05" "09:02" "215A_EL PATO" 28 "LP1912"
Set-Row $ws 24 "08:13:38" "09:41" "215C_EL PATO" 88 "LP1912"
Set-Row $ws 25 "08:34:05" "09:42" "215C_EL PATO" 68 "LP1912"
Set-Row $ws 26 "08:34:05" "10:27" "215A_EL PATO" 113 "LP1912"

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2, 1).Value = "Última actualización: 08:34:05"
$ws.Cells.Item(3, 1).Value = "Total filas: 24"

Set-Row $ws 26 "08:34:05" "08:46" "215A_LA PLATA" 12 "L6173"
Set-Row $ws 27 "08:13:38" "09:08" "215D_LA PLATA" 55 "L6203"
Set-Row $ws 28 "08:34:05" "09:09" "215D_LA PLATA" 35 "L6203"
Set-Row $ws 29 "08:34:05" "10:03" "215B_LP-P MOR-40 Y 115" 89 "L6173"
